# Insert two new columns (O and P) before the old "Extracted Objects" column.
# This shifts old columns O:U to Q:W.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1:P1").EntireColumn.Insert()

# Rename the headers of the (now unchanged-letter) columns M and N.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Set headers for the two newly inserted columns O and P.
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Populate the new O/P column data for rows 2-6.
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1

$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 4

$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 2

$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 4

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
